# Add a new "Skill Description" column (full skill name) right after
# the existing "SkillCode" column, shifting SFIA Level / Keycode /
# Description one column to the right (B->C, C->D, D->E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column B; everything from old B onward shifts right.
$ws.Columns.Item(2).Insert()

# Map each short SkillCode to its full descriptive name. Codes that are
# already full words (Autonomy, Influence, Complexity, Knowledge) map to
# themselves.
$nameMap = @{
    'INVA' = 'Investment appraisal'
    'FMIT' = 'Financial management'
    'PEMT' = 'Performance management'
}

# Header row.
$ws.Cells.Item(1, 2).Value2 = 'Skill Description'

# Find the last used row (column A still holds the SkillCode values).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($code)) {
        continue
    }
    if ($nameMap.ContainsKey($code)) {
        $fullName = $nameMap[$code]
    } else {
        $fullName = $code
    }
    $ws.Cells.Item($r, 2).Value2 = $fullName
}
